# Update the marksheet's correct/total marks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: number of right answers used for marking (row 11)
$ws.Range("B11").Value = 5

# "Total" row: total marks scored out of max, both the raw score cell
# and the "scored/max" text summary (row 12)
$ws.Range("B12").Value = 110
$ws.Range("E12").Value = "110/140"
